$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held 8 data rows (rows 2-9): one row per
# Sending-cluster x Target-cluster(MuSCs/ECs) combination, computed from an
# older TPM dataset. The refreshed script now emits exactly one row per
# Sending cluster (Target cluster is always MuSCs) with updated numbers, so
# the trailing three rows (7-9) are removed entirely and rows 2-6 are
# rewritten in place with the new values.

$ws.Rows("7:9").Delete()

# Row 2: ECs -> Cxcl13/Ccr10 -> MuSCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2087556666666667
$ws.Range("H2").Value = 0.626267
$ws.Range("I2").Value = 0.004212716022507852
$ws.Range("J2").Value = 0.004233737959305679
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.1452775
$ws.Range("N2").Value = 0.290555
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.03032750136416667
$ws.Range("R2").Value = 0.181965008185
$ws.Range("S2").Value = 0.004212716022507852
$ws.Range("T2").Value = 0.004233737959305679

# Row 3: FAPs -> Cxcl13/Ccr10 -> MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 48.06441766666666
$ws.Range("H3").Value = 144.193253
$ws.Range("I3").Value = 0.9699460888896084
$ws.Range("J3").Value = 0.9747862312749473
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.1452775
$ws.Range("N3").Value = 0.290555
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 6.982678437569167
$ws.Range("R3").Value = 41.896070625415
$ws.Range("S3").Value = 0.9699460888896084
$ws.Range("T3").Value = 0.9747862312749473

# Row 4: Inflammatory-Mac -> Cxcl13/Ccr10 -> MuSCs
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.245141
$ws.Range("H4").Value = 0.735423
$ws.Range("I4").Value = 0.004946976697512072
$ws.Range("J4").Value = 0.004971662679410635
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.1452775
$ws.Range("N4").Value = 0.290555
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.0356134716275
$ws.Range("R4").Value = 0.213680829765
$ws.Range("S4").Value = 0.004946976697512072
$ws.Range("T4").Value = 0.004971662679410635

# Row 5: MuSCs -> Cxcl13/Ccr10 -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.7381525
$ws.Range("H5").Value = 1.476305
$ws.Range("I5").Value = 0.01489601175123818
$ws.Range("J5").Value = 0.009980229707158081
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.1452775
$ws.Range("N5").Value = 0.290555
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.10723694981875
$ws.Range("R5").Value = 0.428947799275
$ws.Range("S5").Value = 0.01489601175123818
$ws.Range("T5").Value = 0.009980229707158081

# Row 6: Resolving-Mac -> Cxcl13/Ccr10 -> MuSCs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2972333333333333
$ws.Range("H6").Value = 0.8917
$ws.Range("I6").Value = 0.005998206639133552
$ws.Range("J6").Value = 0.006028138379178328
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.1452775
$ws.Range("N6").Value = 0.290555
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0.04318131558333334
$ws.Range("R6").Value = 0.2590878935
$ws.Range("S6").Value = 0.005998206639133552
$ws.Range("T6").Value = 0.006028138379178328
